# Refresh churn report with the 7 Oct 2024 data pull.
# The header (row 1) and existing company rows (2-273) are unchanged;
# this appends the 71 newly-churned companies as rows 274-344.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 274
$ws.Cells.Item(274, 1).Value = "A/S MOGENS FREDERIKSEN AUTOMOBILER"
$ws.Cells.Item(274, 2).NumberFormat = "@"
$ws.Cells.Item(274, 2).Value = "11"
$ws.Cells.Item(274, 2).Style = "Normal"
$ws.Cells.Item(274, 3).Value = 0

# Row 275
$ws.Cells.Item(275, 1).Value = "RENGØRINGSSELSKABET DRIFT APS"
$ws.Cells.Item(275, 2).Value = 0
$ws.Cells.Item(275, 3).Value = 0

# Row 276
$ws.Cells.Item(276, 1).Value = "SAMSØE & SAMSØE WHOLE SALE APS"
$ws.Cells.Item(276, 2).NumberFormat = "@"
$ws.Cells.Item(276, 2).Value = "1"
$ws.Cells.Item(276, 2).Style = "Normal"
$ws.Cells.Item(276, 3).Value = 226

# Row 277
$ws.Cells.Item(277, 1).Value = "AJ VACCINES A/S"
$ws.Cells.Item(277, 2).NumberFormat = "@"
$ws.Cells.Item(277, 2).Value = "30"
$ws.Cells.Item(277, 2).Style = "Normal"
$ws.Cells.Item(277, 3).Value = 45

# Row 278
$ws.Cells.Item(278, 1).Value = "MICKI LEKSAKER AB"
$ws.Cells.Item(278, 2).Value = 0
$ws.Cells.Item(278, 3).Value = 65

# Row 279
$ws.Cells.Item(279, 1).Value = "COPENHAGEN INFRASTRUCTURE SERVICE COMPANY APS"
$ws.Cells.Item(279, 2).NumberFormat = "@"
$ws.Cells.Item(279, 2).Value = "1"
$ws.Cells.Item(279, 2).Style = "Normal"
$ws.Cells.Item(279, 3).Value = 76

# Row 280
$ws.Cells.Item(280, 1).Value = "TAURUS EJENDOMSADMINISTRATION APS"
$ws.Cells.Item(280, 2).NumberFormat = "@"
$ws.Cells.Item(280, 2).Value = "1"
$ws.Cells.Item(280, 2).Style = "Normal"
$ws.Cells.Item(280, 3).Value = 109

# Row 281
$ws.Cells.Item(281, 1).Value = "NNIT A/S"
$ws.Cells.Item(281, 2).Value = 0
$ws.Cells.Item(281, 3).Value = 117

# Row 282
$ws.Cells.Item(282, 1).Value = "PROTECTOR FORSIKRING DANMARK, FILIAL AF PROTECTOR FORSIKRING ASA, NORGE"
$ws.Cells.Item(282, 2).NumberFormat = "@"
$ws.Cells.Item(282, 2).Value = "4"
$ws.Cells.Item(282, 2).Style = "Normal"
$ws.Cells.Item(282, 3).Value = 127

# Row 283
$ws.Cells.Item(283, 1).Value = "REJSEKORT & REJSEPLAN A/S"
$ws.Cells.Item(283, 2).NumberFormat = "@"
$ws.Cells.Item(283, 2).Value = "46"
$ws.Cells.Item(283, 2).Style = "Normal"
$ws.Cells.Item(283, 3).Value = 0

# Row 284
$ws.Cells.Item(284, 1).Value = "REJSEKORT & REJSEPLAN A/S"
$ws.Cells.Item(284, 2).NumberFormat = "@"
$ws.Cells.Item(284, 2).Value = "46"
$ws.Cells.Item(284, 2).Style = "Normal"
$ws.Cells.Item(284, 3).Value = 0

# Row 285
$ws.Cells.Item(285, 1).Value = "SANTANDER CONSUMER BANK, FILIAL AF SANTANDER CONSUMER BANK AS, NORGE"
$ws.Cells.Item(285, 2).NumberFormat = "@"
$ws.Cells.Item(285, 2).Value = "1"
$ws.Cells.Item(285, 2).Style = "Normal"
$ws.Cells.Item(285, 3).Value = 296

# Row 286
$ws.Cells.Item(286, 1).Value = "SANTANDER CONSUMER BANK, FILIAL AF SANTANDER CONSUMER BANK AS, NORGE"
$ws.Cells.Item(286, 2).NumberFormat = "@"
$ws.Cells.Item(286, 2).Value = "1"
$ws.Cells.Item(286, 2).Style = "Normal"
$ws.Cells.Item(286, 3).Value = 296

# Row 287
$ws.Cells.Item(287, 1).Value = "GULDAGER A/S"
$ws.Cells.Item(287, 2).Value = 0
$ws.Cells.Item(287, 3).Value = 70

# Row 288
$ws.Cells.Item(288, 1).Value = "BP SOENDERGAARD A/S"
$ws.Cells.Item(288, 2).NumberFormat = "@"
$ws.Cells.Item(288, 2).Value = "2"
$ws.Cells.Item(288, 2).Style = "Normal"
$ws.Cells.Item(288, 3).Value = 121

# Row 289
$ws.Cells.Item(289, 1).Value = "RVV A.M.B.A."
$ws.Cells.Item(289, 2).NumberFormat = "@"
$ws.Cells.Item(289, 2).Value = "4"
$ws.Cells.Item(289, 2).Style = "Normal"
$ws.Cells.Item(289, 3).Value = 0

# Row 290
$ws.Cells.Item(290, 1).Value = "ASAPACK A/S"
$ws.Cells.Item(290, 2).NumberFormat = "@"
$ws.Cells.Item(290, 2).Value = "19"
$ws.Cells.Item(290, 2).Style = "Normal"
$ws.Cells.Item(290, 3).Value = 0

# Row 291
$ws.Cells.Item(291, 1).Value = "ROCHE INNOVATION CENTER COPENHAGEN A/S"
$ws.Cells.Item(291, 2).Value = 0
$ws.Cells.Item(291, 3).Value = 134

# Row 292
$ws.Cells.Item(292, 1).Value = "MUJI DENMARK APS"
$ws.Cells.Item(292, 2).Value = 0
$ws.Cells.Item(292, 3).Value = 0

# Row 293
$ws.Cells.Item(293, 1).Value = "KECON A/S"
$ws.Cells.Item(293, 2).NumberFormat = "@"
$ws.Cells.Item(293, 2).Value = "2"
$ws.Cells.Item(293, 2).Style = "Normal"
$ws.Cells.Item(293, 3).Value = 0

# Row 294
$ws.Cells.Item(294, 1).Value = "BAM DANMARK A/S"
$ws.Cells.Item(294, 2).NumberFormat = "@"
$ws.Cells.Item(294, 2).Value = "12"
$ws.Cells.Item(294, 2).Style = "Normal"
$ws.Cells.Item(294, 3).Value = 0

# Row 295
$ws.Cells.Item(295, 1).Value = "FONDET FOR DANSK-NORSK SAMARBEJDE SCHÆFFERGÅRDEN"
$ws.Cells.Item(295, 2).NumberFormat = "@"
$ws.Cells.Item(295, 2).Value = "2"
$ws.Cells.Item(295, 2).Style = "Normal"
$ws.Cells.Item(295, 3).Value = 0

# Row 296
$ws.Cells.Item(296, 1).Value = "AFRY APS"
$ws.Cells.Item(296, 2).NumberFormat = "@"
$ws.Cells.Item(296, 2).Value = "44"
$ws.Cells.Item(296, 2).Style = "Normal"
$ws.Cells.Item(296, 3).Value = 0

# Row 297
$ws.Cells.Item(297, 2).Value = 0
$ws.Cells.Item(297, 3).Value = 0

# Row 298
$ws.Cells.Item(298, 1).Value = "BRANDE BUSLINIER APS"
$ws.Cells.Item(298, 2).NumberFormat = "@"
$ws.Cells.Item(298, 2).Value = "76"
$ws.Cells.Item(298, 2).Style = "Normal"
$ws.Cells.Item(298, 3).Value = 0

# Row 299
$ws.Cells.Item(299, 1).Value = "BRANDE BUSLINIER APS"
$ws.Cells.Item(299, 2).NumberFormat = "@"
$ws.Cells.Item(299, 2).Value = "76"
$ws.Cells.Item(299, 2).Style = "Normal"
$ws.Cells.Item(299, 3).Value = 0

# Row 300
$ws.Cells.Item(300, 1).Value = "ORGANON DENMARK APS"
$ws.Cells.Item(300, 2).Value = 0
$ws.Cells.Item(300, 3).Value = 62

# Row 301
$ws.Cells.Item(301, 1).Value = "2020 SUPERIOR A/S"
$ws.Cells.Item(301, 2).Value = 0
$ws.Cells.Item(301, 3).Value = 0

# Row 302
$ws.Cells.Item(302, 2).Value = 0
$ws.Cells.Item(302, 3).Value = 0

# Row 303
$ws.Cells.Item(303, 1).Value = "NINE UNITED SCANDINAVIA A/S"
$ws.Cells.Item(303, 2).NumberFormat = "@"
$ws.Cells.Item(303, 2).Value = "3"
$ws.Cells.Item(303, 2).Style = "Normal"
$ws.Cells.Item(303, 3).Value = 0

# Row 304
$ws.Cells.Item(304, 1).Value = "NINE UNITED SCANDINAVIA A/S"
$ws.Cells.Item(304, 2).NumberFormat = "@"
$ws.Cells.Item(304, 2).Value = "3"
$ws.Cells.Item(304, 2).Style = "Normal"
$ws.Cells.Item(304, 3).Value = 0

# Row 305
$ws.Cells.Item(305, 1).Value = "NINE UNITED SCANDINAVIA A/S"
$ws.Cells.Item(305, 2).NumberFormat = "@"
$ws.Cells.Item(305, 2).Value = "3"
$ws.Cells.Item(305, 2).Style = "Normal"
$ws.Cells.Item(305, 3).Value = 0

# Row 306
$ws.Cells.Item(306, 1).Value = "IPTOR DENMARK A/S"
$ws.Cells.Item(306, 2).NumberFormat = "@"
$ws.Cells.Item(306, 2).Value = "2"
$ws.Cells.Item(306, 2).Style = "Normal"
$ws.Cells.Item(306, 3).Value = 0

# Row 307
$ws.Cells.Item(307, 1).Value = "NIELSEN CAR GROUP A/S"
$ws.Cells.Item(307, 2).NumberFormat = "@"
$ws.Cells.Item(307, 2).Value = "83"
$ws.Cells.Item(307, 2).Style = "Normal"
$ws.Cells.Item(307, 3).Value = 0

# Row 308
$ws.Cells.Item(308, 1).Value = "OJ RÅDGIVENDE INGENIØRER A/S"
$ws.Cells.Item(308, 2).NumberFormat = "@"
$ws.Cells.Item(308, 2).Value = "22"
$ws.Cells.Item(308, 2).Style = "Normal"
$ws.Cells.Item(308, 3).Value = 0

# Row 309
$ws.Cells.Item(309, 1).Value = "ALMAS KORN A/S"
$ws.Cells.Item(309, 2).NumberFormat = "@"
$ws.Cells.Item(309, 2).Value = "15"
$ws.Cells.Item(309, 2).Style = "Normal"
$ws.Cells.Item(309, 3).Value = 0

# Row 310
$ws.Cells.Item(310, 2).Value = 0
$ws.Cells.Item(310, 3).Value = 0

# Row 311
$ws.Cells.Item(311, 1).Value = "BIANCO FOOTWEAR A/S"
$ws.Cells.Item(311, 2).NumberFormat = "@"
$ws.Cells.Item(311, 2).Value = "12"
$ws.Cells.Item(311, 2).Style = "Normal"
$ws.Cells.Item(311, 3).Value = 0

# Row 312
$ws.Cells.Item(312, 1).Value = "BIANCO FOOTWEAR A/S"
$ws.Cells.Item(312, 2).NumberFormat = "@"
$ws.Cells.Item(312, 2).Value = "12"
$ws.Cells.Item(312, 2).Style = "Normal"
$ws.Cells.Item(312, 3).Value = 0

# Row 313
$ws.Cells.Item(313, 1).Value = "NINE UNITED SCANDINAVIA A/S"
$ws.Cells.Item(313, 2).NumberFormat = "@"
$ws.Cells.Item(313, 2).Value = "3"
$ws.Cells.Item(313, 2).Style = "Normal"
$ws.Cells.Item(313, 3).Value = 0

# Row 314
$ws.Cells.Item(314, 1).Value = "NINE UNITED SCANDINAVIA A/S"
$ws.Cells.Item(314, 2).NumberFormat = "@"
$ws.Cells.Item(314, 2).Value = "3"
$ws.Cells.Item(314, 2).Style = "Normal"
$ws.Cells.Item(314, 3).Value = 0

# Row 315
$ws.Cells.Item(315, 1).Value = "NINE UNITED SCANDINAVIA A/S"
$ws.Cells.Item(315, 2).NumberFormat = "@"
$ws.Cells.Item(315, 2).Value = "3"
$ws.Cells.Item(315, 2).Style = "Normal"
$ws.Cells.Item(315, 3).Value = 0

# Row 316
$ws.Cells.Item(316, 1).Value = "VISBJERGGÅRDEN A/S"
$ws.Cells.Item(316, 2).Value = 0
$ws.Cells.Item(316, 3).Value = 0

# Row 317
$ws.Cells.Item(317, 1).Value = "MSC DENMARK A/S"
$ws.Cells.Item(317, 2).Value = 0
$ws.Cells.Item(317, 3).Value = 0

# Row 318
$ws.Cells.Item(318, 1).Value = "BECH-BRUUN ADVOKATPARTNERSELSKAB"
$ws.Cells.Item(318, 2).NumberFormat = "@"
$ws.Cells.Item(318, 2).Value = "4"
$ws.Cells.Item(318, 2).Style = "Normal"
$ws.Cells.Item(318, 3).Value = 0

# Row 319
$ws.Cells.Item(319, 1).Value = "GROUP ONLINE A/S"
$ws.Cells.Item(319, 2).NumberFormat = "@"
$ws.Cells.Item(319, 2).Value = "18"
$ws.Cells.Item(319, 2).Style = "Normal"
$ws.Cells.Item(319, 3).Value = 0

# Row 320
$ws.Cells.Item(320, 1).Value = "BLUE HORS APS"
$ws.Cells.Item(320, 2).Value = 0
$ws.Cells.Item(320, 3).Value = 0

# Row 321
$ws.Cells.Item(321, 1).Value = "NORTH MEDIA A/S"
$ws.Cells.Item(321, 2).Value = 0
$ws.Cells.Item(321, 3).Value = 0

# Row 322
$ws.Cells.Item(322, 2).Value = 0
$ws.Cells.Item(322, 3).Value = 0

# Row 323
$ws.Cells.Item(323, 1).Value = "EF EDUCATION A/S"
$ws.Cells.Item(323, 2).NumberFormat = "@"
$ws.Cells.Item(323, 2).Value = "1"
$ws.Cells.Item(323, 2).Style = "Normal"
$ws.Cells.Item(323, 3).Value = 85

# Row 324
$ws.Cells.Item(324, 1).Value = "LEICA GEOSYSTEMS TECHNOLOGY A/S"
$ws.Cells.Item(324, 2).NumberFormat = "@"
$ws.Cells.Item(324, 2).Value = "45"
$ws.Cells.Item(324, 2).Style = "Normal"
$ws.Cells.Item(324, 3).Value = 0

# Row 325
$ws.Cells.Item(325, 1).Value = "DANSK ENERGI MANAGEMENT A/S"
$ws.Cells.Item(325, 2).NumberFormat = "@"
$ws.Cells.Item(325, 2).Value = "17"
$ws.Cells.Item(325, 2).Style = "Normal"
$ws.Cells.Item(325, 3).Value = 0

# Row 326
$ws.Cells.Item(326, 1).Value = "DANSK ENERGI MANAGEMENT A/S"
$ws.Cells.Item(326, 2).NumberFormat = "@"
$ws.Cells.Item(326, 2).Value = "17"
$ws.Cells.Item(326, 2).Style = "Normal"
$ws.Cells.Item(326, 3).Value = 0

# Row 327
$ws.Cells.Item(327, 1).Value = "DANSK ENERGI MANAGEMENT A/S"
$ws.Cells.Item(327, 2).NumberFormat = "@"
$ws.Cells.Item(327, 2).Value = "17"
$ws.Cells.Item(327, 2).Style = "Normal"
$ws.Cells.Item(327, 3).Value = 0

# Row 328
$ws.Cells.Item(328, 1).Value = "DANSK ENERGI MANAGEMENT A/S"
$ws.Cells.Item(328, 2).NumberFormat = "@"
$ws.Cells.Item(328, 2).Value = "17"
$ws.Cells.Item(328, 2).Style = "Normal"
$ws.Cells.Item(328, 3).Value = 0

# Row 329
$ws.Cells.Item(329, 1).Value = "IVECO DANMARK A/S"
$ws.Cells.Item(329, 2).NumberFormat = "@"
$ws.Cells.Item(329, 2).Value = "1"
$ws.Cells.Item(329, 2).Style = "Normal"
$ws.Cells.Item(329, 3).Value = 126

# Row 330
$ws.Cells.Item(330, 1).Value = "ESPRIT DE CORP. DANMARK A/S"
$ws.Cells.Item(330, 2).NumberFormat = "@"
$ws.Cells.Item(330, 2).Value = "11"
$ws.Cells.Item(330, 2).Style = "Normal"
$ws.Cells.Item(330, 3).Value = 0

# Row 331
$ws.Cells.Item(331, 1).Value = "HENNING IBSEN A/S"
$ws.Cells.Item(331, 2).NumberFormat = "@"
$ws.Cells.Item(331, 2).Value = "13"
$ws.Cells.Item(331, 2).Style = "Normal"
$ws.Cells.Item(331, 3).Value = 0

# Row 332
$ws.Cells.Item(332, 1).Value = "ADVENT TECHNOLOGIES A/S"
$ws.Cells.Item(332, 2).NumberFormat = "@"
$ws.Cells.Item(332, 2).Value = "12"
$ws.Cells.Item(332, 2).Style = "Normal"
$ws.Cells.Item(332, 3).Value = 0

# Row 333
$ws.Cells.Item(333, 1).Value = "SABRO A/S"
$ws.Cells.Item(333, 2).NumberFormat = "@"
$ws.Cells.Item(333, 2).Value = "2"
$ws.Cells.Item(333, 2).Style = "Normal"
$ws.Cells.Item(333, 3).Value = 0

# Row 334
$ws.Cells.Item(334, 1).Value = "KYNDRYL DANMARK APS"
$ws.Cells.Item(334, 2).Value = 0
$ws.Cells.Item(334, 3).Value = 155

# Row 335
$ws.Cells.Item(335, 1).Value = "PHOENIX X A/S"
$ws.Cells.Item(335, 2).NumberFormat = "@"
$ws.Cells.Item(335, 2).Value = "1"
$ws.Cells.Item(335, 2).Style = "Normal"
$ws.Cells.Item(335, 3).Value = 57

# Row 336
$ws.Cells.Item(336, 1).Value = "PHOENIX X A/S"
$ws.Cells.Item(336, 2).NumberFormat = "@"
$ws.Cells.Item(336, 2).Value = "1"
$ws.Cells.Item(336, 2).Style = "Normal"
$ws.Cells.Item(336, 3).Value = 57

# Row 337
$ws.Cells.Item(337, 1).Value = "EPINION P/S"
$ws.Cells.Item(337, 2).NumberFormat = "@"
$ws.Cells.Item(337, 2).Value = "1"
$ws.Cells.Item(337, 2).Style = "Normal"
$ws.Cells.Item(337, 3).Value = 146

# Row 338
$ws.Cells.Item(338, 1).Value = "COPENHAGEN INFRASTRUCTURE SERVICE COMPANY APS"
$ws.Cells.Item(338, 2).NumberFormat = "@"
$ws.Cells.Item(338, 2).Value = "1"
$ws.Cells.Item(338, 2).Style = "Normal"
$ws.Cells.Item(338, 3).Value = 76

# Row 339
$ws.Cells.Item(339, 1).Value = "GEOPARTNER LANDINSPEKTØRER A/S"
$ws.Cells.Item(339, 2).NumberFormat = "@"
$ws.Cells.Item(339, 2).Value = "15"
$ws.Cells.Item(339, 2).Style = "Normal"
$ws.Cells.Item(339, 3).Value = 0

# Row 340
$ws.Cells.Item(340, 1).Value = "LØNMODTAGERNES DYRTIDSFOND"
$ws.Cells.Item(340, 2).NumberFormat = "@"
$ws.Cells.Item(340, 2).Value = "1"
$ws.Cells.Item(340, 2).Style = "Normal"
$ws.Cells.Item(340, 3).Value = 18

# Row 341
$ws.Cells.Item(341, 1).Value = "LØNMODTAGERNES DYRTIDSFOND"
$ws.Cells.Item(341, 2).NumberFormat = "@"
$ws.Cells.Item(341, 2).Value = "1"
$ws.Cells.Item(341, 2).Style = "Normal"
$ws.Cells.Item(341, 3).Value = 18

# Row 342
$ws.Cells.Item(342, 1).Value = "JESPERS TORVEKØKKEN APS"
$ws.Cells.Item(342, 2).NumberFormat = "@"
$ws.Cells.Item(342, 2).Value = "30"
$ws.Cells.Item(342, 2).Style = "Normal"
$ws.Cells.Item(342, 3).Value = 0

# Row 343
$ws.Cells.Item(343, 1).Value = "BLUJAY SOLUTIONS A/S"
$ws.Cells.Item(343, 2).NumberFormat = "@"
$ws.Cells.Item(343, 2).Value = "32"
$ws.Cells.Item(343, 2).Style = "Normal"
$ws.Cells.Item(343, 3).Value = 0

# Row 344
$ws.Cells.Item(344, 1).Value = "SOL OG STRAND FERIEHUSUDLEJNING A/S"
$ws.Cells.Item(344, 2).NumberFormat = "@"
$ws.Cells.Item(344, 2).Value = "12"
$ws.Cells.Item(344, 2).Style = "Normal"
$ws.Cells.Item(344, 3).Value = 0
